# Update the build version/timestamp references throughout the workbook.
$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# Row 2: "Version: ..." on the About sheet
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# Row 6: Recommended citation text on the About sheet
$aboutSheet.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Saranskaya Coal Mine, Kazakhstan, M1436, version ' + "'" + $newVersion + "'" + '. (See the CC license for attribution requirements if sharing or adapting the data set.)'

# Column S (build_version) for each data row on the "Boundaries and methane sources" sheet
$usedRange = $dataSheet.UsedRange
$lastRow = $usedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $dataSheet.Range("S" + $r)
    $current = $cell.Text
    if ($current -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
